$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 3.755628166162433),
    @(0.2881169905109251, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 1, 1.276001064745008),
    @(0.6545652718822623, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 0, 1.642449346116345),
    @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 0, 17.65757632934944),
    @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 3.755628166162433),
    @(0.1169995834814548, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 1, 1.104883657715537),
    @(0.1169995834814548, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 1, 1.104883657715537)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}
